$wb = $excel.ActiveWorkbook

# Rename the second sheet from "fdcryvy" to "sadgdda"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "sadgdda"

# Update recalculated values (slope comparison of previous/current R peak used to discard T waves)
$ws = $ws2
$ws.Cells.Item(1, 10).Value = 37.0259473323822
$ws.Cells.Item(2, 10).Value = 48.22730660438538
$ws.Cells.Item(3, 10).Value = 43.76771211624146
$ws.Cells.Item(4, 2).Value = 2593
$ws.Cells.Item(4, 4).Value = 2561
$ws.Cells.Item(4, 5).Value = 31
$ws.Cells.Item(4, 6).Value = 10
$ws.Cells.Item(4, 7).Value = 99.61104628549202
$ws.Cells.Item(4, 8).Value = 98.80401234567901
$ws.Cells.Item(4, 9).Value = 0.01594090202177294
$ws.Cells.Item(4, 10).Value = 37.52178311347961
$ws.Cells.Item(5, 10).Value = 44.612628698349
$ws.Cells.Item(6, 10).Value = 38.75931739807129
$ws.Cells.Item(7, 10).Value = 43.0322003364563
$ws.Cells.Item(8, 10).Value = 37.83224892616272
$ws.Cells.Item(9, 10).Value = 42.17259883880615
$ws.Cells.Item(10, 2).Value = 1810
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 8).Value = 99.11553344389165
$ws.Cells.Item(10, 9).Value = 0.009470752089136491
$ws.Cells.Item(10, 10).Value = 43.17305111885071
$ws.Cells.Item(11, 10).Value = 34.08266234397888
$ws.Cells.Item(12, 10).Value = 37.57371211051941
$ws.Cells.Item(13, 10).Value = 36.49975299835205
$ws.Cells.Item(14, 10).Value = 36.02704548835754
$ws.Cells.Item(15, 10).Value = 37.84151339530945
$ws.Cells.Item(16, 10).Value = 40.8201630115509
$ws.Cells.Item(17, 10).Value = 42.01306819915771
$ws.Cells.Item(18, 10).Value = 42.40998435020447
$ws.Cells.Item(19, 10).Value = 41.25351238250732
$ws.Cells.Item(20, 10).Value = 32.56008791923523
$ws.Cells.Item(21, 2).Value = 2600
$ws.Cells.Item(21, 4).Value = 2595
$ws.Cells.Item(21, 6).Value = 5
$ws.Cells.Item(21, 7).Value = 99.80769230769231
$ws.Cells.Item(21, 8).Value = 99.84609465178914
$ws.Cells.Item(21, 9).Value = 0.003460207612456748
$ws.Cells.Item(21, 10).Value = 38.84509062767029
$ws.Cells.Item(22, 10).Value = 45.72270345687866
$ws.Cells.Item(23, 2).Value = 2128
$ws.Cells.Item(23, 4).Value = 2127
$ws.Cells.Item(23, 6).Value = 8
$ws.Cells.Item(23, 7).Value = 99.62529274004685
$ws.Cells.Item(23, 9).Value = 0.003745318352059925
$ws.Cells.Item(23, 10).Value = 34.53421521186829
$ws.Cells.Item(24, 2).Value = 2939
$ws.Cells.Item(24, 4).Value = 2917
$ws.Cells.Item(24, 5).Value = 21
$ws.Cells.Item(24, 6).Value = 62
$ws.Cells.Item(24, 7).Value = 97.91876468613629
$ws.Cells.Item(24, 8).Value = 99.28522804628999
$ws.Cells.Item(24, 9).Value = 0.02785234899328859
$ws.Cells.Item(24, 10).Value = 40.53076672554016
$ws.Cells.Item(25, 2).Value = 2644
$ws.Cells.Item(25, 4).Value = 2643
$ws.Cells.Item(25, 6).Value = 12
$ws.Cells.Item(25, 7).Value = 99.54802259887006
$ws.Cells.Item(25, 9).Value = 0.004518072289156626
$ws.Cells.Item(25, 10).Value = 40.32425355911255
$ws.Cells.Item(26, 10).Value = 34.93050694465637
$ws.Cells.Item(27, 10).Value = 39.12341213226318
$ws.Cells.Item(28, 2).Value = 2957
$ws.Cells.Item(28, 4).Value = 2954
$ws.Cells.Item(28, 6).Value = 50
$ws.Cells.Item(28, 7).Value = 98.33555259653795
$ws.Cells.Item(28, 8).Value = 99.93234100135318
$ws.Cells.Item(28, 9).Value = 0.0173044925124792
$ws.Cells.Item(28, 10).Value = 38.8621723651886
$ws.Cells.Item(29, 2).Value = 2643
$ws.Cells.Item(29, 4).Value = 2635
$ws.Cells.Item(29, 5).Value = 7
$ws.Cells.Item(29, 6).Value = 14
$ws.Cells.Item(29, 7).Value = 99.4714986787467
$ws.Cells.Item(29, 8).Value = 99.73504920514762
$ws.Cells.Item(29, 9).Value = 0.007924528301886792
$ws.Cells.Item(29, 10).Value = 41.02306056022644
$ws.Cells.Item(30, 10).Value = 40.03092861175537
$ws.Cells.Item(31, 10).Value = 39.33634877204895
$ws.Cells.Item(32, 2).Value = 2260
$ws.Cells.Item(32, 5).Value = 1
$ws.Cells.Item(32, 8).Value = 99.95573262505533
$ws.Cells.Item(32, 9).Value = 0.001768346595932803
$ws.Cells.Item(32, 10).Value = 41.80965399742126
$ws.Cells.Item(33, 2).Value = 3361
$ws.Cells.Item(33, 4).Value = 3360
$ws.Cells.Item(33, 6).Value = 2
$ws.Cells.Item(33, 7).Value = 99.94051160023795
$ws.Cells.Item(33, 9).Value = 0.0005947071067499256
$ws.Cells.Item(33, 10).Value = 38.53718090057373
$ws.Cells.Item(34, 10).Value = 42.34335994720459
$ws.Cells.Item(35, 2).Value = 2047
$ws.Cells.Item(35, 4).Value = 2046
$ws.Cells.Item(35, 6).Value = 1
$ws.Cells.Item(35, 7).Value = 99.95114802149487
$ws.Cells.Item(35, 9).Value = 0.00048828125
$ws.Cells.Item(35, 10).Value = 48.87689232826233
$ws.Cells.Item(36, 10).Value = 39.58051323890686
$ws.Cells.Item(37, 2).Value = 2462
$ws.Cells.Item(37, 4).Value = 2460
$ws.Cells.Item(37, 6).Value = 22
$ws.Cells.Item(37, 7).Value = 99.11361804995971
$ws.Cells.Item(37, 8).Value = 99.95936611133685
$ws.Cells.Item(37, 9).Value = 0.009262988320579944
$ws.Cells.Item(37, 10).Value = 39.62017297744751
$ws.Cells.Item(38, 2).Value = 2605
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(38, 8).Value = 100
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 34.49544548988342
$ws.Cells.Item(39, 2).Value = 2067
$ws.Cells.Item(39, 5).Value = 20
$ws.Cells.Item(39, 8).Value = 99.03194578896418
$ws.Cells.Item(39, 9).Value = 0.0126643935703848
$ws.Cells.Item(39, 10).Value = 42.82982897758484
$ws.Cells.Item(40, 2).Value = 2256
$ws.Cells.Item(40, 4).Value = 2254
$ws.Cells.Item(40, 5).Value = 1
$ws.Cells.Item(40, 6).Value = 1
$ws.Cells.Item(40, 7).Value = 99.95565410199556
$ws.Cells.Item(40, 8).Value = 99.95565410199556
$ws.Cells.Item(40, 10).Value = 44.55140495300293
$ws.Cells.Item(41, 10).Value = 40.57601618766785
$ws.Cells.Item(42, 10).Value = 39.79496383666992
$ws.Cells.Item(43, 10).Value = 45.0031943321228
$ws.Cells.Item(44, 2).Value = 2752
$ws.Cells.Item(44, 4).Value = 2751
$ws.Cells.Item(44, 6).Value = 1
$ws.Cells.Item(44, 7).Value = 99.96366279069767
$ws.Cells.Item(44, 9).Value = 0.0003632401017072285
$ws.Cells.Item(44, 10).Value = 40.2917058467865

